$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0
$ws.Cells.Item(2, 3).Value = 0.00223048327137546
$ws.Cells.Item(2, 4).Value = 0.0423791821561338
$ws.Cells.Item(2, 5).Value = 0.00223048327137546
$ws.Cells.Item(2, 6).Value = 0.992565055762082
$ws.Cells.Item(2, 7).Value = 0.041635687732342
$ws.Cells.Item(2, 8).Value = 0.00371747211895911
$ws.Cells.Item(2, 9).Value = 0.888475836431227
$ws.Cells.Item(2, 10).Value = 0
$ws.Cells.Item(2, 11).Value = 0
$ws.Cells.Item(2, 12).Value = 0.642379182156134
$ws.Cells.Item(2, 13).Value = 0.00223048327137546
$ws.Cells.Item(2, 14).Value = 0.980669144981413
$ws.Cells.Item(2, 15).Value = 0.985873605947955
$ws.Cells.Item(2, 16).Value = 0.179182156133829
$ws.Cells.Item(2, 17).Value = 0.568029739776952
$ws.Cells.Item(2, 18).Value = 0.990334572490706
$ws.Cells.Item(2, 19).Value = 0.00297397769516729
$ws.Cells.Item(2, 20).Value = 0.00297397769516729
$ws.Cells.Item(2, 21).Value = 0.0111524163568773
$ws.Cells.Item(2, 22).Value = 0.150185873605948
$ws.Cells.Item(2, 23).Value = 0.00520446096654275
$ws.Cells.Item(2, 24).Value = 0.0111524163568773

$ws.Cells.Item(3, 2).Value = 0.993308550185874
$ws.Cells.Item(3, 3).Value = 0.89814126394052
$ws.Cells.Item(3, 4).Value = 0.00520446096654275
$ws.Cells.Item(3, 5).Value = 0.00594795539033457
$ws.Cells.Item(3, 6).Value = 0.000743494423791822
$ws.Cells.Item(3, 7).Value = 0.0148698884758364
$ws.Cells.Item(3, 8).Value = 0
$ws.Cells.Item(3, 9).Value = 0
$ws.Cells.Item(3, 10).Value = 0.987360594795539
$ws.Cells.Item(3, 11).Value = 0.00148698884758364
$ws.Cells.Item(3, 12).Value = 0.000743494423791822
$ws.Cells.Item(3, 13).Value = 0.950185873605948
$ws.Cells.Item(3, 14).Value = 0.00446096654275093
$ws.Cells.Item(3, 15).Value = 0.00446096654275093
$ws.Cells.Item(3, 16).Value = 0.00817843866171004
$ws.Cells.Item(3, 17).Value = 0.00223048327137546
$ws.Cells.Item(3, 18).Value = 0.00520446096654275
$ws.Cells.Item(3, 19).Value = 0.863940520446097
$ws.Cells.Item(3, 20).Value = 0.00817843866171004
$ws.Cells.Item(3, 21).Value = 0
$ws.Cells.Item(3, 22).Value = 0.739776951672863
$ws.Cells.Item(3, 23).Value = 0.00297397769516729
$ws.Cells.Item(3, 24).Value = 0.0163568773234201

$ws.Cells.Item(4, 2).Value = 0.00148698884758364
$ws.Cells.Item(4, 3).Value = 0.00594795539033457
$ws.Cells.Item(4, 4).Value = 0.850557620817844
$ws.Cells.Item(4, 5).Value = 0.991078066914498
$ws.Cells.Item(4, 6).Value = 0.00297397769516729
$ws.Cells.Item(4, 7).Value = 0.937546468401487
$ws.Cells.Item(4, 8).Value = 0.00297397769516729
$ws.Cells.Item(4, 9).Value = 0.0921933085501859
$ws.Cells.Item(4, 10).Value = 0.000743494423791822
$ws.Cells.Item(4, 11).Value = 0
$ws.Cells.Item(4, 12).Value = 0.00594795539033457
$ws.Cells.Item(4, 13).Value = 0.00371747211895911
$ws.Cells.Item(4, 14).Value = 0.0118959107806691
$ws.Cells.Item(4, 15).Value = 0.00446096654275093
$ws.Cells.Item(4, 16).Value = 0.811895910780669
$ws.Cells.Item(4, 17).Value = 0.331598513011152
$ws.Cells.Item(4, 18).Value = 0.00446096654275093
$ws.Cells.Item(4, 19).Value = 0.0966542750929368
$ws.Cells.Item(4, 20).Value = 0.00297397769516729
$ws.Cells.Item(4, 21).Value = 0.984386617100372
$ws.Cells.Item(4, 22).Value = 0.00892193308550186
$ws.Cells.Item(4, 23).Value = 0.107063197026022
$ws.Cells.Item(4, 24).Value = 0.883271375464684

$ws.Cells.Item(5, 2).Value = 0.00520446096654275
$ws.Cells.Item(5, 3).Value = 0.0936802973977695
$ws.Cells.Item(5, 4).Value = 0.10185873605948
$ws.Cells.Item(5, 5).Value = 0
$ws.Cells.Item(5, 6).Value = 0.00371747211895911
$ws.Cells.Item(5, 7).Value = 0.00594795539033457
$ws.Cells.Item(5, 8).Value = 0.993308550185874
$ws.Cells.Item(5, 9).Value = 0.0193308550185874
$ws.Cells.Item(5, 10).Value = 0.0118959107806691
$ws.Cells.Item(5, 11).Value = 0.998513011152416
$ws.Cells.Item(5, 12).Value = 0.35092936802974
$ws.Cells.Item(5, 13).Value = 0.0438661710037175
$ws.Cells.Item(5, 14).Value = 0.00297397769516729
$ws.Cells.Item(5, 15).Value = 0.00520446096654275
$ws.Cells.Item(5, 16).Value = 0.000743494423791822
$ws.Cells.Item(5, 17).Value = 0.0981412639405204
$ws.Cells.Item(5, 18).Value = 0
$ws.Cells.Item(5, 19).Value = 0.0356877323420074
$ws.Cells.Item(5, 20).Value = 0.985873605947955
$ws.Cells.Item(5, 21).Value = 0.00446096654275093
$ws.Cells.Item(5, 22).Value = 0.101115241635688
$ws.Cells.Item(5, 23).Value = 0.884758364312268
$ws.Cells.Item(5, 24).Value = 0.0892193308550186

